$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.534.33"
$ws.Range("E2").Value = "  +3.15%  "
$ws.Range("D3").Value = "1.692.99"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.42%  "
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("E7").Value = "  +1.43%  "
$ws.Range("E8").Value = "  +1.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.524"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.001"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.16"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08723"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.203"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.82%  "
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.566"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.83%  "
$ws.Range("D17").Value = "1.690.65"
$ws.Range("E17").Value = "  +0.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "99.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07047"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.859"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.09%  "
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("E23").Value = "  +1.19%  "
$ws.Range("D24").Value = "24.533.18"
$ws.Range("E24").Value = "  +3.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.069"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.325"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.88"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.220"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.540"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.78%  "
$ws.Range("D32").Value = "1.878.81"
$ws.Range("E32").Value = "  +1.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.086"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08532"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.274"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.25"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.10%  "
$ws.Range("E37").Value = "  -1.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2701"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.39"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02737"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09002"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.469"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7639"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7149"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.517"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.196"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "140.49"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.325"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07987"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.70%  "
